$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert new column B (SEX) -- shifts old B:X to C:Y, carrying values + styles
$ws.Columns("B:B").Insert()

# 2. Insert 3 new rows above the (now row 2) data row, pushing it down to row 5
$ws.Rows("2:4").Insert()

# 3. Make sure new header cell Y1 matches the bold/border/centered style of the other headers
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)

# 4. Force the data rows to literal text so ISO-looking dates are not auto-converted to serials
$ws.Range("A2:Y5").NumberFormat = "@"

# 5. Row 1 headers -> UPPER_CASE / snake_case; add SEX + OTHER_COMMENTS
$ws.Cells.Item(1,2).Value = "SEX"
$ws.Cells.Item(1,3).Value = "NAME"
$ws.Cells.Item(1,4).Value = "TEAM"
$ws.Cells.Item(1,5).Value = "CODE"
$ws.Cells.Item(1,6).Value = "INJURY_DATE"
$ws.Cells.Item(1,7).Value = "RETURN_DATE"
$ws.Cells.Item(1,8).Value = "INJURY_LOCATION"
$ws.Cells.Item(1,9).Value = "INJURY_SIDE"
$ws.Cells.Item(1,10).Value = "INJURY_TYPE"
$ws.Cells.Item(1,11).Value = "OCCURRENCE"
$ws.Cells.Item(1,12).Value = "OVERUSE_TRAUMA"
$ws.Cells.Item(1,13).Value = "ONSET"
$ws.Cells.Item(1,14).Value = "CONTACT"
$ws.Cells.Item(1,15).Value = "ACTION"
$ws.Cells.Item(1,16).Value = "ACTION_DESCRIPTION"
$ws.Cells.Item(1,17).Value = "RE_INJURY"
$ws.Cells.Item(1,18).Value = "REFEREE_SANCTION"
$ws.Cells.Item(1,19).Value = "DIAGNOSTIC_EXAMINATION"
$ws.Cells.Item(1,20).Value = "DIAGNOSIS"
$ws.Cells.Item(1,21).Value = "SURGERY"
$ws.Cells.Item(1,22).Value = "MENSTRUAL_PHASE"
$ws.Cells.Item(1,23).Value = "ORAL_CONTRACEPTIVES"
$ws.Cells.Item(1,24).Value = "HORMONAL_CONTRACEPTIVES"
$ws.Cells.Item(1,25).Value = "OTHER_COMMENTS"

# 6. Row 2 (new record: 1.docx / Male, rest blank)
$ws.Cells.Item(2,1).Value = "1.docx"
$ws.Cells.Item(2,2).Value = "Male"

# 7. Row 3 (new record: 2.docx / Male)
$ws.Cells.Item(3,1).Value = "2.docx"
$ws.Cells.Item(3,2).Value = "Male"
$ws.Cells.Item(3,3).Value = "Luca F"
$ws.Cells.Item(3,4).Value = "FC B"
$ws.Cells.Item(3,5).Value = "xx11x"
$ws.Cells.Item(3,6).Value = "2025-06-23"
$ws.Cells.Item(3,7).Value = "Wrong date format"
$ws.Cells.Item(3,8).Value = "Lower leg (incl. Achilles tendon)"
$ws.Cells.Item(3,9).Value = "Bilateral/central"
$ws.Cells.Item(3,10).Value = "Muscle rupture/tear/strain**"
$ws.Cells.Item(3,11).Value = "N/A (gradual onset injury)"
$ws.Cells.Item(3,12).Value = "Overuse (repetitive mechanism)"
$ws.Cells.Item(3,13).Value = "Gradual onset"
$ws.Cells.Item(3,14).Value = "Yes"
$ws.Cells.Item(3,17).Value = "No"
$ws.Cells.Item(3,18).Value = "No foul, Opponent foul, Own foul, Yellow card, Red card"
$ws.Cells.Item(3,19).Value = "Ultrasonography"
$ws.Cells.Item(3,21).Value = "No"
$ws.Cells.Item(3,23).Value = "Yes"
$ws.Cells.Item(3,24).Value = "Should appear"
$ws.Cells.Item(3,25).Value = "None"

# 8. Row 4 (new record: 3.docx / Female)
$ws.Cells.Item(4,1).Value = "3.docx"
$ws.Cells.Item(4,2).Value = "Female"
$ws.Cells.Item(4,3).Value = "All filled"
$ws.Cells.Item(4,4).Value = "FC ALL"
$ws.Cells.Item(4,5).Value = "XXXXX"
$ws.Cells.Item(4,6).Value = "2025-06-10"
$ws.Cells.Item(4,7).Value = "2025-10-12"
$ws.Cells.Item(4,8).Value = "Head, Abdomen, Elbow, Hip, Lower leg (incl. Achilles tendon), Neck, Lumbosacral, Forearm, Groin, Ankle, Chest, Shoulder, Wrist, Thigh, Foot, Thoracic spine, Upper arm, Hand, Knee"
$ws.Cells.Item(4,9).Value = "Right, Left, Bilateral/central"
$ws.Cells.Item(4,10).Value = "Concussion, Meniscus lesion, Haematoma/contusion/bruise (incl. compartment syndrome), Fracture (specify if stress fracture), Cartilage lesion, Nerve injury (central or peripheral other than concussion), Other bone injury (e.g., bone stress), Muscle rupture/tear/strain**, Dental injury*, Joint dislocation/subluxation*, Tendon rupture/tendinopathy, Vessel injury (excl. skin haematoma), Joint sprain (i.e., ligament/capsule), Abrasion, Bursitis, Arthritis/synovitis/capsulitis, Laceration, Overuse unspecified, Additional Injury"
$ws.Cells.Item(4,11).Value = "Training, Match (min. of injury: 70), N/A (gradual onset injury), Football training, Football & other training, League match, Friendly match, Other training, Reserve/youth team training, Champions League, Reserve/youth team match, National team, Other cup match"
$ws.Cells.Item(4,12).Value = "Too many answers"
$ws.Cells.Item(4,13).Value = "Too many answers"
$ws.Cells.Item(4,14).Value = "Too many answers"
$ws.Cells.Item(4,15).Value = "Running/sprinting, Heading, Controlling the ball, Tackling other player, Blocked*, Twisting/turning, Landing (incl. jumping), Hit by ball, Tackled by other player, Use of arm/elbow*, Shooting/passing, Falling (incl. diving), Collision, Sliding/stretching*, Other player action"
$ws.Cells.Item(4,16).Value = "Did everything"
$ws.Cells.Item(4,17).Value = "Too many answers"
$ws.Cells.Item(4,18).Value = "No foul, Opponent foul, Own foul, Yellow card, Red card"
$ws.Cells.Item(4,19).Value = "Clinical only, X-ray, Ultrasonography, MRI, Other"
$ws.Cells.Item(4,20).Value = "DIagnosis"
$ws.Cells.Item(4,21).Value = "Too many answers"
$ws.Cells.Item(4,22).Value = "Too many answers"
$ws.Cells.Item(4,23).Value = "Too many answers"
$ws.Cells.Item(4,24).Value = "Too many answers"
$ws.Cells.Item(4,25).Value = "Comment"

# 9. Row 5 (original record, now shifted) -- fix SEX + ISO dates + simplified Contact value
$ws.Cells.Item(5,2).Value = "Female"
$ws.Cells.Item(5,6).Value = "2022-05-03"
$ws.Cells.Item(5,7).Value = "2022-06-03"
$ws.Cells.Item(5,14).Value = "Yes"
